$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B15").Value = 11000025
$ws.Range("C15").Value = 8000025
$ws.Range("C17").Value = 8000069
$ws.Range("C18").Value = 8000075
$ws.Range("C20").Value = 8000034
$ws.Range("D20").Value = "https://www.areu.lombardia.it/web/home/soreu-dei-laghi"
$ws.Range("E20").Value = "SOREU dei Laghi"
$ws.Range("C21").Value = 8000032
$ws.Range("D21").Value = "http://www.protezionecivile.gov.it/servizio-nazionale/strutture-operative/volontariato"
$ws.Range("E21").Value = $null
$ws.Range("B22").Value = 11000037
$ws.Range("C22").Value = 8000037
$ws.Range("D22").Value = "https://www.areu.lombardia.it/web/home/soreu"
$ws.Range("E22").Value = "Le SOREU"
$ws.Range("B23").Value = 11000038
$ws.Range("C23").Value = 8000038
$ws.Range("D23").Value = "https://www4.ti.ch/di/smpp/chi-siamo/presentazione/"
$ws.Range("E23").Value = "Chi siamo"
$ws.Range("B24").Value = 11000039
$ws.Range("B25").Value = 11000040
$ws.Range("C25").Value = 8000040
$ws.Range("B26").Value = 11000041
$ws.Range("C26").Value = 8000041
$ws.Range("B27").Value = 11000043
$ws.Range("C27").Value = 8000043
$ws.Range("B28").Value = 11000044
$ws.Range("C28").Value = 8000044
$ws.Range("B29").Value = 11000045
$ws.Range("C29").Value = 8000068
$ws.Range("D29").Value = "http://www.protezionecivile.gov.it/dipartimento"
$ws.Range("E29").Value = "Dipartimento"
$ws.Range("B30").Value = 11000046
$ws.Range("C30").Value = 8000046
$ws.Range("D30").Value = "http://www.irpi.cnr.it/focus/suscettibilita-da-frana/"
$ws.Range("E30").Value = "Modelli e carte di suscettibilità da frana"
$ws.Range("B31").Value = 11000047
$ws.Range("C31").Value = 8000047
$ws.Range("E31").Value = "Legge federale sulla protezione della popolazione e sulla protezione civile del 4 ottobre 2002"
$ws.Range("B32").Value = 11000047
$ws.Range("C32").Value = 8000069
$ws.Range("D32").Value = "https://www.admin.ch/opc/it/classified-compilation/20011872/201701010000/520.1.pdf"
$ws.Range("E32").Value = "-- documento confronto normativa -- wp 3.2 gestisco -- da completare"
$ws.Range("B33").Value = 11000048
$ws.Range("C33").Value = 8000019
$ws.Range("D33").Value = "https://www.gazzettaufficiale.it/eli/id/2018/1/22/18G00011/sg"
$ws.Range("E33").Value = "Decreto Legislativo 2 Gennaio 2018, N. 1, Codice Della Protezione Civile. (18G00011)"
$ws.Range("B34").Value = 11000053
$ws.Range("C34").Value = 8000053
$ws.Range("D34").Value = "https://m3.ti.ch/CAN/RLeggi/public/index.php/raccolta-leggi/legge/num/48"
$ws.Range("E34").Value = "Legge sulla protezione della popolazione (del 26 febbraio 2007)"
$ws.Range("B35").Value = 11000057
$ws.Range("C35").Value = 8000057
$ws.Range("D35").Value = "https://www.babs.admin.ch/content/babs-internet/it/publikservice/downloads/unterlagen-ausbildung/_jcr_content/contentPar/accordion_1920886228/accordionItems/kommando_zivilschutz/accordionPar/downloadlist_copy/downloadItems/829_1459931125997.download/personal170191103it.pdf"
$ws.Range("E35").Value = "Il comando della protezione civile - Personale"
$ws.Range("B36").Value = 11000059
$ws.Range("C36").Value = 8000059
$ws.Range("D36").Value = "https://m3.ti.ch/CAN/RLeggi/public/index.php/index/nuovafinestra/atto/49/volume/5%20SICUREZZA/numLegge/500.110"
$ws.Range("E36").Value = "Regolamento sulla protezione della popolazione (RProtPop) (del 18 ottobre 2017)"
$ws.Range("B37").Value = 11000063
$ws.Range("C37").Value = 8000063
$ws.Range("D37").Value = "https://www.ch.ch/it/allarme-sirene/"
$ws.Range("E37").Value = "Segnali di allarme in Svizzera"
$ws.Range("B38").Value = 11000066
$ws.Range("C38").Value = 8000066
$ws.Range("D38").Value = "https://www4.ti.ch/di/smpp/chi-siamo/servizio-della-protezione-della-popolazione/"
$ws.Range("E38").Value = "Servizio della protezione della popolazione"
$ws.Range("B39").Value = 11000073
